$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 24.333402
$ws.Range("H2").Value = 73.000206
$ws.Range("I2").Value = 0.9697175080062574
$ws.Range("J2").Value = 0.9697175080062576
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 45.95651366666667
$ws.Range("N2").Value = 137.869541
$ws.Range("O2").Value = 0.6189188856627118
$ws.Range("P2").Value = 0.6189188856627118
$ws.Range("Q2").Value = 1118.278321569494
$ws.Range("R2").Value = 10064.50489412545
$ws.Range("S2").Value = 0.6001764794628547
$ws.Range("T2").Value = 0.6001764794628547

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 24.333402
$ws.Range("H3").Value = 73.000206
$ws.Range("I3").Value = 0.9697175080062574
$ws.Range("J3").Value = 0.9697175080062576
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 6.849914666666667
$ws.Range("N3").Value = 20.549744
$ws.Range("O3").Value = 0.09225115688993263
$ws.Range("P3").Value = 0.09225115688993261
$ws.Range("Q3").Value = 166.681727249696
$ws.Range("R3").Value = 1500.135545247264
$ws.Range("S3").Value = 0.08945756196999975
$ws.Range("T3").Value = 0.08945756196999975

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 24.333402
$ws.Range("H4").Value = 73.000206
$ws.Range("I4").Value = 0.9697175080062574
$ws.Range("J4").Value = 0.9697175080062576
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 21.446458
$ws.Range("N4").Value = 64.33937399999999
$ws.Range("O4").Value = 0.2888299574473556
$ws.Range("P4").Value = 0.2888299574473556
$ws.Range("Q4").Value = 521.8652839901159
$ws.Range("R4").Value = 4696.787555911043
$ws.Range("S4").Value = 0.2800834665734031
$ws.Range("T4").Value = 0.2800834665734031

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.288887
$ws.Range("H5").Value = 0.866661
$ws.Range("I5").Value = 0.01151252018667195
$ws.Range("J5").Value = 0.01151252018667195
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 45.95651366666667
$ws.Range("N5").Value = 137.869541
$ws.Range("O5").Value = 0.6189188856627118
$ws.Range("P5").Value = 0.6189188856627118
$ws.Range("Q5").Value = 13.27623936362233
$ws.Range("R5").Value = 119.486154272601
$ws.Range("S5").Value = 0.007125316165104479
$ws.Range("T5").Value = 0.007125316165104481

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.288887
$ws.Range("H6").Value = 0.866661
$ws.Range("I6").Value = 0.01151252018667195
$ws.Range("J6").Value = 0.01151252018667195
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 6.849914666666667
$ws.Range("N6").Value = 20.549744
$ws.Range("O6").Value = 0.09225115688993263
$ws.Range("P6").Value = 0.09225115688993261
$ws.Range("Q6").Value = 1.978851298309333
$ws.Range("R6").Value = 17.809661684784
$ws.Range("S6").Value = 0.001062043305939191
$ws.Range("T6").Value = 0.001062043305939191

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.288887
$ws.Range("H7").Value = 0.866661
$ws.Range("I7").Value = 0.01151252018667195
$ws.Range("J7").Value = 0.01151252018667195
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 21.446458
$ws.Range("N7").Value = 64.33937399999999
$ws.Range("O7").Value = 0.2888299574473556
$ws.Range("P7").Value = 0.2888299574473556
$ws.Range("Q7").Value = 6.195602912245999
$ws.Range("R7").Value = 55.76042621021399
$ws.Range("S7").Value = 0.003325160715628283
$ws.Range("T7").Value = 0.003325160715628283

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.4710003333333333
$ws.Range("H8").Value = 1.413001
$ws.Range("I8").Value = 0.01876997180707065
$ws.Range("J8").Value = 0.01876997180707065
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 45.95651366666667
$ws.Range("N8").Value = 137.869541
$ws.Range("O8").Value = 0.6189188856627118
$ws.Range("P8").Value = 0.6189188856627118
$ws.Range("Q8").Value = 21.64553325583789
$ws.Range("R8").Value = 194.809799302541
$ws.Range("S8").Value = 0.01161709003475268
$ws.Range("T8").Value = 0.01161709003475268

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.4710003333333333
$ws.Range("H9").Value = 1.413001
$ws.Range("I9").Value = 0.01876997180707065
$ws.Range("J9").Value = 0.01876997180707065
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.849914666666667
$ws.Range("N9").Value = 20.549744
$ws.Range("O9").Value = 0.09225115688993263
$ws.Range("P9").Value = 0.09225115688993261
$ws.Range("Q9").Value = 3.226312091304889
$ws.Range("R9").Value = 29.036808821744
$ws.Range("S9").Value = 0.001731551613993686
$ws.Range("T9").Value = 0.001731551613993686

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.4710003333333333
$ws.Range("H10").Value = 1.413001
$ws.Range("I10").Value = 0.01876997180707065
$ws.Range("J10").Value = 0.01876997180707065
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 21.446458
$ws.Range("N10").Value = 64.33937399999999
$ws.Range("O10").Value = 0.2888299574473556
$ws.Range("P10").Value = 0.2888299574473556
$ws.Range("Q10").Value = 10.10128886681933
$ws.Range("R10").Value = 90.91159980137398
$ws.Range("S10").Value = 0.005421330158324279
$ws.Range("T10").Value = 0.005421330158324279
